# fix auth bug, remove extra sheet from demo
#
# 1. Correct the title text in A1 on the DemoFunctions sheet (was missing
#    "Boardflare" and had an incorrect capitalization).
# 2. Highlight the title row (A1:C1) with the green accent fill.
# 3. Remove the unused blank "Sheet1" worksheet that was left over in the
#    demo workbook.
# 4. Leave the selection on C13 (the last formula cell) as the active cell.

$wb = $excel.ActiveWorkbook

$demo = $wb.Worksheets("DemoFunctions")

# --- Fix the title text in A1 -------------------------------------------
$demo.Range("A1").Value = "Example functions built using Boardflare Python for Excel"

# --- Highlight the title row with the green accent theme color ----------
$titleRow = $demo.Range("A1:C1")
$titleRow.Interior.ThemeColor = 10

# --- Remove the extra, empty demo sheet ----------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets("Sheet1").Delete()
$excel.DisplayAlerts = $true

# --- Restore selection on the DemoFunctions sheet ------------------------
$demo.Activate()
$demo.Range("C13").Select()
